$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows above the existing data (at row 2), pushing the
# current data (rows 2-21) down to rows 10-29.
$ws.Rows("2:9").Insert()

# The inserted rows inherit the header row's style; clear that back to
# the default (unstyled) so the new data rows match the plain numeric
# cell formatting used by the rest of the table.
$ws.Rows("2:9").ClearFormats()

# New rows to insert at the top (rows 2-9)
$topRows = @(
  @(0.037449836730957, 0.3890565633773803, -0.2103400230407714),
  @(0.0105371475219726, 0.3668201565742492, -0.3450851440429687),
  @(0.1156568527221679, 0.2621434330940246, 0.2320724725723266),
  @(-0.012700080871582, 0.3827533721923828, 0.0918664336204528),
  @(-0.1776895523071289, 1.103561997413635, -1.161585211753845),
  @(-0.1967945098876953, -0.0881298780441284, -0.4354097247123718),
  @(0.3726930618286133, 0.3928739428520202, -0.1955753564834594),
  @(0.1368236541748047, 0.4095092415809631, -0.2153286337852478)
)

$r = 2
foreach ($row in $topRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# New rows appended at the bottom (rows 30-31)
$bottomRows = @(
  @(1.180892944335938, -0.3624088764190674, 1.944910764694214),
  @(-0.6099348068237305, -0.0995303392410278, 1.559979677200317)
)

$r = 30
foreach ($row in $bottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
